$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (deleted/adjusted passive ROM submax values)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Match the new selection range reflected in the saved workbook
$ws.Range("B1:E3").Select()
